# Updated cryptos list on Sat Jun  3 19:22:25 UTC 2023 with GitHub Actions
#
# Writes new Price (column D) and Volume(1h) (column E) figures, and for the
# three rows whose ranking swapped with a neighbour also rewrites the Coin
# name (B) and Link (C).
#
# Column D holds values that look numeric ("1.002", "10.10", "0.000008600",
# "27.207.37", ...). Plain Range.Value assignment lets Excel's automatic
# type-inference turn these into IEEE doubles, which silently drops
# significant trailing zeros (e.g. "10.10" -> 10.1) or collapses multi-dot
# strings. Forcing the cell to Text format before the write, then restoring
# the "Normal" style afterwards (so no stray formatting sticks around),
# guarantees the literal text is preserved exactly as authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

function Set-Text($cellRef, $text) {
    $ws.Range($cellRef).Value = $text
}

# Row 2 - Bitcoin
Set-PriceText "D2" "27.207.37"
Set-Text "E2" "  +0.13%  "

# Row 3 - Ethereum
Set-PriceText "D3" "1.895.25"
Set-Text "E3" "  -0.69%  "

# Row 4 - TetherUSD
Set-PriceText "D4" "1.003"
Set-Text "E4" "  +0.28%  "

# Row 5 - BNB
Set-PriceText "D5" "307.22"
Set-Text "E5" "  -0.18%  "

# Row 6 - USDC
Set-Text "E6" "  +0.22%  "

# Row 7 - XRP
Set-PriceText "D7" "0.5228"
Set-Text "E7" "  -0.37%  "

# Row 8 - Cardano
Set-PriceText "D8" "0.3753"
Set-Text "E8" "  -0.71%  "

# Row 9 - Dogecoin
Set-PriceText "D9" "0.07271"
Set-Text "E9" "  -0.11%  "

# Row 10 - Solana
Set-Text "E10" "  -0.50%  "

# Row 11 - Polygon
Set-PriceText "D11" "0.8991"
Set-Text "E11" "  +0.05%  "

# Row 12 - TRON
Set-PriceText "D12" "0.08162"
Set-Text "E12" "  +6.33%  "

# Row 13 - Litecoin
Set-PriceText "D13" "96.86"
Set-Text "E13" "  +2.22%  "

# Row 14 - WrappedEther
Set-PriceText "D14" "1.897.60"
Set-Text "E14" "  +0.58%  "

# Row 15 - Polkadot
Set-PriceText "D15" "5.275"
Set-Text "E15" "  +0.45%  "

# Row 16 - BinanceUSD
Set-PriceText "D16" "1.003"
Set-Text "E16" "  +0.24%  "

# Row 17 - ShibaInu
Set-PriceText "D17" "0.000008600"
Set-Text "E17" "  +0.59%  "

# Row 19 - Dai
Set-Text "E19" "  +0.21%  "

# Row 20 - WrappedBTC
Set-PriceText "D20" "27.239.65"
Set-Text "E20" "  +0.02%  "

# Row 21 - Uniswap
Set-PriceText "D21" "5.083"
Set-Text "E21" "  -0.13%  "

# Row 22 - Cosmos
Set-PriceText "D22" "10.69"
Set-Text "E22" "  +0.50%  "

# Row 23 - Chainlink
Set-PriceText "D23" "6.400"
Set-Text "E23" "  -0.72%  "

# Row 24 - Monero
Set-PriceText "D24" "147.71"
Set-Text "E24" "  +1.27%  "

# Row 25 - LidoDAOToken
Set-PriceText "D25" "2.290"
Set-Text "E25" "  -1.18%  "

# Row 26 - was EthereumClassic, now Toncoin (rows 26/27 swapped order)
Set-Text "B26" "Toncoin"
Set-Text "C26" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-PriceText "D26" "1.747"
Set-Text "E26" "  +1.02%  "

# Row 27 - was Toncoin, now EthereumClassic
Set-Text "B27" "EthereumClassic"
Set-Text "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-PriceText "D27" "18.22"
Set-Text "E27" "  +0.29%  "

# Row 28 - BitcoinCash
Set-PriceText "D28" "115.06"
Set-Text "E28" "  +0.25%  "

# Row 29 - was InternetComputer(DFINITY), now Filecoin (rows 29/30 swapped order)
Set-Text "B29" "Filecoin"
Set-Text "C29" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-PriceText "D29" "4.914"
Set-Text "E29" "  -1.06%  "

# Row 30 - was Filecoin, now InternetComputer(DFINITY)
Set-Text "B30" "InternetComputer(DFINITY)"
Set-Text "C30" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-PriceText "D30" "4.805"
Set-Text "E30" "  -0.15%  "

# Row 31 - Stellar
Set-PriceText "D31" "0.09224"
Set-Text "E31" "  +0.11%  "

# Row 32 - Hedera
Set-PriceText "D32" "0.05043"
Set-Text "E32" "  -0.69%  "

# Row 33 - ImmutableX
Set-PriceText "D33" "0.7941"
Set-Text "E33" "  +1.90%  "

# Row 34 - ARBITRUM
Set-PriceText "D34" "1.217"
Set-Text "E34" "  -2.33%  "

# Row 35 - MXToken
Set-PriceText "D35" "3.439"
Set-Text "E35" "  +3.92%  "

# Row 36 - HuobiToken
Set-PriceText "D36" "2.946"
Set-Text "E36" "  -1.57%  "

# Row 37 - RenderToken
Set-PriceText "D37" "2.564"
Set-Text "E37" "  -1.70%  "

# Row 38 - TheSandbox
Set-PriceText "D38" "0.5653"
Set-Text "E38" "  -0.50%  "

# Row 39 - VeChain
Set-PriceText "D39" "0.01983"
Set-Text "E39" "  -0.64%  "

# Row 40 - TrustWalletToken
Set-Text "E40" "  +0.14%  "

# Row 41 - Aptos
Set-PriceText "D41" "8.935"
Set-Text "E41" "  -1.04%  "

# Row 42 - FraxShare
Set-PriceText "D42" "6.534"
Set-Text "E42" "  -1.48%  "

# Row 43 - Quant
Set-PriceText "D43" "115.10"
Set-Text "E43" "  -3.01%  "

# Row 44 - Algorand
Set-Text "E44" "  -0.34%  "

# Row 45 - Decentraland
Set-PriceText "D45" "0.4879"
Set-Text "E45" "  +0.48%  "

# Row 46 - PaxDollar
Set-Text "E46" "  +0.20%  "

# Row 47 - EnergySwap
Set-PriceText "D47" "10.10"
Set-Text "E47" "  -1.01%  "

# Row 48 - NEARProtocol
Set-PriceText "D48" "1.617"
Set-Text "E48" "  +0.64%  "

# Row 50 - Aave
Set-Text "E50" "  -1.37%  "

# Row 51 - Cronos
Set-PriceText "D51" "0.05941"
Set-Text "E51" "  +0.21%  "
